# cs1501_rec1_Sept6.pptx -- "fixed typo on slide 6"
#
# Two small text fixes were made to the deck:
#   1. Slide 2 ("Contact info"): the "Note" line was missing an opening
#      parenthesis -- "...different TA see course website..." should read
#      "...different TA (see course website...".
#   2. Slide 6 ("Homework submission through Github"): a typo, "Github
#      will firsts send you an email..." should read "...will first
#      send you an email...".

$p = $ppt.ActivePresentation

# --- Slide 2: add the missing "(" before "see course website..." ---
$slide2 = $p.Slides.Item(2)
$noteShape = $slide2.Shapes.Item(2)
$noteRange = $noteShape.TextFrame.TextRange

$seeIdx = $noteRange.Text.IndexOf("see course")
if ($seeIdx -ge 0) {
    # "see " is 4 characters; prepend "(" to turn it into "(see "
    $seeRun = $noteRange.Characters($seeIdx + 1, 4)
    $seeRun.Text = "(" + $seeRun.Text
}

# --- Slide 6: fix "firsts" -> "first" ---
$slide6 = $p.Slides.Item(6)
$ghShape = $slide6.Shapes.Item(2)
$ghRange = $ghShape.TextFrame.TextRange

$typoIdx = $ghRange.Text.IndexOf("firsts ")
if ($typoIdx -ge 0) {
    $typoRun = $ghRange.Characters($typoIdx + 1, 7)
    $typoRun.Text = "first "
}
